$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.804
$ws.Range("H2").Value = 0.407
$ws.Range("G3").Value = 0.8080000000000001
$ws.Range("L3").Value = 0.447
$ws.Range("D4").Value = 0.479
$ws.Range("G4").Value = 0.734
$ws.Range("G5").Value = 0.789
$ws.Range("J6").Value = 0.446
$ws.Range("C7").Value = 0.679
$ws.Range("J7").Value = 0.408
$ws.Range("C8").Value = 0.657
$ws.Range("L8").Value = 0.471
$ws.Range("J9").Value = 0.367
$ws.Range("E10").Value = 0.595
$ws.Range("C11").Value = 0.6909999999999999
$ws.Range("J12").Value = 0.338
$ws.Range("G13").Value = 0.759
$ws.Range("C14").Value = 0.677
$ws.Range("L14").Value = 0.455
$ws.Range("G15").Value = 0.756
$ws.Range("G16").Value = 0.742
$ws.Range("L16").Value = 0.43
$ws.Range("J17").Value = 0.325
$ws.Range("H18").Value = 0.36
$ws.Range("D20").Value = 0.473
$ws.Range("G21").Value = 0.8080000000000001
$ws.Range("G22").Value = 0.727
$ws.Range("I22").Value = 0.549
$ws.Range("I23").Value = 0.581
$ws.Range("L23").Value = 0.446
$ws.Range("G24").Value = 0.704
$ws.Range("C25").Value = 0.621
$ws.Range("E25").Value = 0.554
$ws.Range("C28").Value = 0.631
$ws.Range("E30").Value = 0.5590000000000001
$ws.Range("E31").Value = 0.548
$ws.Range("F31").Value = 0.503
$ws.Range("C32").Value = 0.571
$ws.Range("H32").Value = 0.463
$ws.Range("H35").Value = 0.374
$ws.Range("L37").Value = 0.432
$ws.Range("H38").Value = 0.401
$ws.Range("L38").Value = 0.46
$ws.Range("F43").Value = 0.521
$ws.Range("G43").Value = 0.626
$ws.Range("G44").Value = 0.72
$ws.Range("H45").Value = 0.372
$ws.Range("C46").Value = 0.628
$ws.Range("J46").Value = 0.413
$ws.Range("D47").Value = 0.502
$ws.Range("L47").Value = 0.453
$ws.Range("F48").Value = 0.554
$ws.Range("G48").Value = 0.671
$ws.Range("E49").Value = 0.533
$ws.Range("I49").Value = 0.602
$ws.Range("L52").Value = 0.432
$ws.Range("G54").Value = 0.669
$ws.Range("L54").Value = 0.444
$ws.Range("E56").Value = 0.547
$ws.Range("J56").Value = 0.431
$ws.Range("G57").Value = 0.6899999999999999
$ws.Range("H57").Value = 0.367
$ws.Range("L58").Value = 0.431
$ws.Range("G59").Value = 0.676
$ws.Range("L59").Value = 0.446
$ws.Range("E60").Value = 0.538
$ws.Range("L61").Value = 0.433
$ws.Range("G62").Value = 0.599
$ws.Range("D63").Value = 0.444
$ws.Range("G63").Value = 0.642
$ws.Range("J65").Value = 0.379
$ws.Range("C66").Value = 0.606
$ws.Range("J68").Value = 0.474
$ws.Range("H69").Value = 0.444
$ws.Range("L70").Value = 0.415
$ws.Range("E71").Value = 0.578
$ws.Range("I71").Value = 0.5629999999999999
$ws.Range("L73").Value = 0.407
$ws.Range("L75").Value = 0.415
$ws.Range("C76").Value = 0.662
$ws.Range("J76").Value = 0.404
$ws.Range("E77").Value = 0.511
$ws.Range("B78").Value = 0.591
$ws.Range("H79").Value = 0.441
$ws.Range("I79").Value = 0.607
$ws.Range("E80").Value = 0.536
$ws.Range("G81").Value = 0.649
$ws.Range("H81").Value = 0.385
$ws.Range("G82").Value = 0.772
$ws.Range("K82").Value = 0.392
$ws.Range("C83").Value = 0.618
$ws.Range("F84").Value = 0.653
$ws.Range("H84").Value = 0.475
$ws.Range("K85").Value = 0.441
$ws.Range("E86").Value = 0.523
$ws.Range("I86").Value = 0.501
$ws.Range("D87").Value = 0.416
$ws.Range("E87").Value = 0.545
$ws.Range("B88").Value = 0.705
$ws.Range("J88").Value = 0.471
$ws.Range("C89").Value = 0.638
$ws.Range("J89").Value = 0.362
$ws.Range("H90").Value = 0.447
$ws.Range("F91").Value = 0.5679999999999999
$ws.Range("G91").Value = 0.6899999999999999
$ws.Range("C93").Value = 0.641
$ws.Range("G93").Value = 0.647
$ws.Range("G95").Value = 0.754
$ws.Range("G96").Value = 0.72
$ws.Range("H97").Value = 0.47
$ws.Range("D98").Value = 0.48
$ws.Range("H98").Value = 0.414
$ws.Range("J99").Value = 0.393
$ws.Range("J100").Value = 0.474
